$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "struggle"
$ws.Cells.Item(22, 3).Value = -3.367526054382324
$ws.Cells.Item(22, 4).Value = 0.915987193584442
$ws.Cells.Item(22, 5).Value = -1.626443386077881
$ws.Cells.Item(22, 6).Value = -1.0144944190979
$ws.Cells.Item(22, 7).Value = -1.167210817337036
$ws.Cells.Item(22, 8).Value = 0.6551529765129089

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "struggle"
$ws.Cells.Item(23, 3).Value = -1.118759155273438
$ws.Cells.Item(23, 4).Value = 2.792432069778442
$ws.Cells.Item(23, 5).Value = -4.963344097137451
$ws.Cells.Item(23, 6).Value = -0.3874412775039673
$ws.Cells.Item(23, 7).Value = -0.7050912380218506
$ws.Cells.Item(23, 8).Value = 0.0914770737290382

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "struggle"
$ws.Cells.Item(24, 3).Value = -1.557756900787354
$ws.Cells.Item(24, 4).Value = -0.5582034587860107
$ws.Cells.Item(24, 5).Value = -0.2619988918304443
$ws.Cells.Item(24, 6).Value = -0.1458440721035003
$ws.Cells.Item(24, 7).Value = -0.3762930035591125
$ws.Cells.Item(24, 8).Value = -0.0704022198915481

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "struggle"
$ws.Cells.Item(25, 3).Value = 1.715949058532715
$ws.Cells.Item(25, 4).Value = -1.576748490333557
$ws.Cells.Item(25, 5).Value = 5.096891403198242
$ws.Cells.Item(25, 6).Value = 0.2157881408929824
$ws.Cells.Item(25, 7).Value = 0.3240640163421631
$ws.Cells.Item(25, 8).Value = 0.0951422601938247

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "struggle"
$ws.Cells.Item(26, 3).Value = -0.8243503570556641
$ws.Cells.Item(26, 4).Value = 0.5943599939346313
$ws.Cells.Item(26, 5).Value = 1.927432060241699
$ws.Cells.Item(26, 6).Value = 0.0403171069920063
$ws.Cells.Item(26, 7).Value = 0.1484402567148208
$ws.Cells.Item(26, 8).Value = -0.0852157026529312

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "struggle"
$ws.Cells.Item(27, 3).Value = -0.2981023788452148
$ws.Cells.Item(27, 4).Value = 1.024843096733093
$ws.Cells.Item(27, 5).Value = 0.8517363667488098
$ws.Cells.Item(27, 6).Value = 0.1014036312699318
$ws.Cells.Item(27, 7).Value = 0.3179553747177124
$ws.Cells.Item(27, 8).Value = 0.0390953756868839

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "struggle"
$ws.Cells.Item(28, 3).Value = 0.4371089935302734
$ws.Cells.Item(28, 4).Value = 0.3337190449237823
$ws.Cells.Item(28, 5).Value = -0.154114544391632
$ws.Cells.Item(28, 6).Value = 0.052381694316864
$ws.Cells.Item(28, 7).Value = 0.1099557429552078
$ws.Cells.Item(28, 8).Value = 0.0681114718317985

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "struggle"
$ws.Cells.Item(29, 3).Value = -0.5641984939575195
$ws.Cells.Item(29, 4).Value = -0.3292053341865539
$ws.Cells.Item(29, 5).Value = -0.326197862625122
$ws.Cells.Item(29, 6).Value = 0.0522289797663688
$ws.Cells.Item(29, 7).Value = -0.4196644127368927
$ws.Cells.Item(29, 8).Value = 0.2273945808410644

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "struggle"
$ws.Cells.Item(30, 3).Value = 0.131052017211914
$ws.Cells.Item(30, 4).Value = 0.5107872486114502
$ws.Cells.Item(30, 5).Value = 0.0540084838867187
$ws.Cells.Item(30, 6).Value = 0.0937678143382072
$ws.Cells.Item(30, 7).Value = -0.1565342247486114
$ws.Cells.Item(30, 8).Value = 0.0675006061792373

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "struggle"
$ws.Cells.Item(31, 3).Value = -0.1777238845825195
$ws.Cells.Item(31, 4).Value = 0.4102384448051452
$ws.Cells.Item(31, 5).Value = 0.1352127194404602
$ws.Cells.Item(31, 6).Value = -0.0591012127697467
$ws.Cells.Item(31, 7).Value = 0.0331394411623477
$ws.Cells.Item(31, 8).Value = 0.0291688162833452
